# Data provider concept in DD framework
# Adds three new rows of test data (Somu/Ravi/Uma) to the RegTestData sheet,
# wires up the Emailaddress/Password hyperlinks for each new row, and moves
# the active selection the way the authored workbook left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegTestData")

# Populate the new data column-by-column so shared strings are interned in
# the same order as the source edit (Somu, Ravi, Uma, then R, S, R).
$ws.Range("A3").Value = "Somu"
$ws.Range("A4").Value = "Ravi"
$ws.Range("A5").Value = "Uma"

$ws.Range("B3").Value = "R"
$ws.Range("B4").Value = "S"
$ws.Range("B5").Value = "R"

$ws.Range("C3").Value = "padmapriya@gmail.com"
$ws.Range("C4").Value = "padmapriya@gmail.com"
$ws.Range("C5").Value = "padmapriya@gmail.com"

$ws.Range("D3").Value = "padma@123"
$ws.Range("D4").Value = "padma@123"
$ws.Range("D5").Value = "padma@123"

# Hyperlink the Emailaddress / Password cells exactly like row 2.
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:padmapriya@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:padmapriya@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:padmapriya@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:padma@123")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:padma@123")
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:padma@123")

# Hyperlinks.Add() stamps its own style; restore the shared hyperlink style
# (same one already used by C2/D2) on each new cell.
$ws.Range("C3").Style = $ws.Range("C2").Style
$ws.Range("C4").Style = $ws.Range("C2").Style
$ws.Range("C5").Style = $ws.Range("C2").Style
$ws.Range("D3").Style = $ws.Range("D2").Style
$ws.Range("D4").Style = $ws.Range("D2").Style
$ws.Range("D5").Style = $ws.Range("D2").Style

# Match the author's final selection.
$ws.Range("B11").Select()
